$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark from its current spot
#    (it currently sits right after "Sam McMillan: 5").
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2. Insert the new clause before the final period of the sentence
#    "...very early stage of development."
# ------------------------------------------------------------------
$find = $d.Content
$find.Find.Execute("very early stage of development.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$periodPos = $find.End - 1
$insertPoint = $d.Range($periodPos, $periodPos)
$insertPoint.InsertBefore(", by completing the movement mechanics and scoring system, as well as getting a physical prototype ready for feedback")

# ------------------------------------------------------------------
# 3. Re-add the "_GoBack" bookmark at the very end of that paragraph
#    (just after the trailing period, before the paragraph mark).
# ------------------------------------------------------------------
$idealPara = $d.Paragraphs(12).Range
$newBookmarkPos = $idealPara.End - 1
$bookmarkRange = $d.Range($newBookmarkPos, $newBookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
